$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "jatin bhalla"
$ws.Range("B4").Value = "jatinbhalla18@gmail.com"
$ws.Range("C4").Value = "Devops engineer"

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:jatinbhalla18@gmail.com")
$ws.Range("B4").Style = "Hyperlink"

$ws.Range("C4").Select()
